# Generate Report for Handback
# Updates the status of the bf2592bb-4a59-4284-90bc-bbeb04756fdc.md entry
# from "Ready for handoff" to "Handed back: in sync with en-US" across the
# Overview / zh-cn / de-de sheets, and refreshes the related handback
# timestamps / clears the stale error detail now that the handback is current.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("H3").Value = "2016-08-31 09:00:51"
$zh.Range("K3").Value = "2016-08-31 09:00:51"
$zh.Range("P3").Value = ""
$zh.Columns.Item(16).ColumnWidth = 13.7470528738839

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("K3").Value = "2016-08-31 09:01:07"
$de.Range("P3").Value = ""
$de.Columns.Item(16).ColumnWidth = 13.7470528738839
